# Insert a new weekly price record for "Poroto verde" (Vega Central Mapocho de
# Santiago) above the current row 722, shifting all subsequent rows down by
# one. This matches the commit "Fruta / hortaliza, semanal" which appends the
# latest weekly observation at the top of this sub-sheet's price history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 722:808 down to 723:809 and insert a blank row at 722.
$ws.Rows.Item(722).Insert()

$newRow = 722
$values = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45124, 13, 100112031, "Poroto verde", "Magnum", "Primera", 70, 22000, 23000, 22500, "`$/malla 25 kilos", "Perú", 900, 25, "Hortaliza")

for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $values[$col - 1]
}
